$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C10 previously held 18; the upstream commit restores/overwrites it with 1.
$ws.Range("C10").Value = 1
